# Auto-generated: apply scheduled Kraken price-runner updates to each leve-profit sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 262.6
$ws.Range("I2").Value = 262.6
$ws.Range("K2").Value = 262.6
$ws.Range("M2").Value = -149.6

$ws.Range("H9").Value = 338.33334
$ws.Range("J9").Value = 95
$ws.Range("L9").Value = 95
$ws.Range("N9").Value = -433

$ws.Range("H43").Value = 2865.3125
$ws.Range("I43").Value = 2107.5
$ws.Range("J43").Value = 3320
$ws.Range("K43").Value = 2107.5
$ws.Range("L43").Value = 3320
$ws.Range("M43").Value = -2038.5
$ws.Range("N43").Value = -3458

$ws.Range("H110").Value = 18500
$ws.Range("J110").Value = 18500
$ws.Range("L110").Value = 18500
$ws.Range("N110").Value = -26680

$ws.Range("H116").Value = 2000
$ws.Range("J116").Value = 2000
$ws.Range("L116").Value = 2000
$ws.Range("N116").Value = -8884

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 50228
$ws.Range("J107").Value = 50228
$ws.Range("L107").Value = 50228
$ws.Range("N107").Value = -57908

$ws.Range("H124").Value = 32500
$ws.Range("J124").Value = 32500
$ws.Range("L124").Value = 32500
$ws.Range("N124").Value = -42320

$ws.Range("H130").Value = 94443
$ws.Range("J130").Value = 94443
$ws.Range("L130").Value = 94443
$ws.Range("N130").Value = -104483

$ws.Range("H135").Value = 45000
$ws.Range("J135").Value = 45000
$ws.Range("L135").Value = 45000
$ws.Range("N135").Value = -55140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H112").Value = 20000
$ws.Range("J112").Value = 20000
$ws.Range("L112").Value = 20000
$ws.Range("N112").Value = -22954

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 260.27274
$ws.Range("I7").Value = 184.85715
$ws.Range("J7").Value = 392.25
$ws.Range("K7").Value = 184.85715
$ws.Range("L7").Value = 392.25
$ws.Range("M7").Value = -71.85714999999999
$ws.Range("N7").Value = -618.25

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H116").Value = 77777
$ws.Range("J116").Value = 77777
$ws.Range("L116").Value = 77777
$ws.Range("N116").Value = -86955

$ws.Range("H117").Value = 95000
$ws.Range("J117").Value = 95000
$ws.Range("L117").Value = 95000
$ws.Range("N117").Value = -104178

$ws.Range("H130").Value = 49750
$ws.Range("J130").Value = 49750
$ws.Range("L130").Value = 49750
$ws.Range("N130").Value = -59790

$ws.Range("H134").Value = 5144.75
$ws.Range("I134").Value = 5144.75
$ws.Range("K134").Value = 15434.25
$ws.Range("M134").Value = -12899.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1166777.6
$ws.Range("I4").Value = 1187625
$ws.Range("K4").Value = 3562875
$ws.Range("M4").Value = -3562763

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 8000
$ws.Range("J119").Value = 8000
$ws.Range("L119").Value = 8000
$ws.Range("N119").Value = -17676

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2358.3333
$ws.Range("J2").Value = 2358.3333
$ws.Range("L2").Value = 2358.3333
$ws.Range("N2").Value = -2582.3333

$ws.Range("H7").Value = 2841
$ws.Range("I7").Value = 2481.1667
$ws.Range("K7").Value = 2481.1667
$ws.Range("M7").Value = -2369.1667

$ws.Range("H46").Value = 1867.6666
$ws.Range("I46").Value = 1900.5
$ws.Range("K46").Value = 1900.5
$ws.Range("M46").Value = -1712.5

$ws.Range("H55").Value = 1440.6666
$ws.Range("I55").Value = 1136.125
$ws.Range("J55").Value = 2049.75
$ws.Range("K55").Value = 1136.125
$ws.Range("L55").Value = 2049.75
$ws.Range("M55").Value = -963.125
$ws.Range("N55").Value = -2395.75

$ws.Range("H61").Value = 2928.0715
$ws.Range("I61").Value = 2675.1
$ws.Range("K61").Value = 2675.1
$ws.Range("M61").Value = -2473.1

$ws.Range("H113").Value = 2928.0715
$ws.Range("I113").Value = 2675.1
$ws.Range("K113").Value = 2675.1
$ws.Range("M113").Value = -505.0999999999999

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H126").Value = 2841
$ws.Range("I126").Value = 2481.1667
$ws.Range("K126").Value = 7443.500100000001
$ws.Range("M126").Value = -4973.500100000001

$ws.Range("H128").Value = 89500
$ws.Range("J128").Value = 89500
$ws.Range("L128").Value = 89500
$ws.Range("N128").Value = -99460

$ws.Range("H135").Value = 290000
$ws.Range("J135").Value = 290000
$ws.Range("L135").Value = 290000
$ws.Range("N135").Value = -300140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 887241.4
$ws.Range("I2").Value = 1133558.9
$ws.Range("J2").Value = 498.4
$ws.Range("K2").Value = 1133558.9
$ws.Range("L2").Value = 498.4
$ws.Range("M2").Value = -1133446.9
$ws.Range("N2").Value = -722.4

$ws.Range("H4").Value = 16594.4
$ws.Range("I4").Value = 21376.158
$ws.Range("K4").Value = 21376.158
$ws.Range("M4").Value = -21263.158

$ws.Range("H81").Value = 40000.6
$ws.Range("I81").Value = 40000.6
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 80001.2
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -78940.2
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 40000.6
$ws.Range("I84").Value = 40000.6
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 400006
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -394702
$ws.Range("N84").ClearContents()

$ws.Range("H129").Value = 99980
$ws.Range("J129").Value = 99980
$ws.Range("L129").Value = 99980
$ws.Range("N129").Value = -109980

$ws.Range("H135").Value = 99700
$ws.Range("J135").Value = 99700
$ws.Range("L135").Value = 99700
$ws.Range("N135").Value = -109840

$ws.Range("H137").Value = 97715
$ws.Range("J137").Value = 97715
$ws.Range("L137").Value = 97715
$ws.Range("N137").Value = -107915

Write-Output "Kraken profit sheets updated."